# Daily attendance processing - 2025-10-16 08:55:01
# Normalizes the "Recorded By" column (G) so that the most recent recorder
# ("System") is listed first: rotate each comma-separated list of
# recorders one place to the right (last entry moves to the front).
# Cells holding only a single recorder are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $value = $cell.Text

    if ([string]::IsNullOrEmpty($value)) {
        continue
    }

    $parts = $value -split ', '

    if ($parts.Count -gt 1) {
        $rotated = @($parts[-1]) + $parts[0..($parts.Count - 2)]
        $cell.Value = [string]::Join(', ', $rotated)
    }
}
